$d = $word.ActiveDocument

# Title paragraph: "Questions: Introduction to complex numbers"
$titleRange = $d.Paragraphs.Item(1).Range
$titleRange.Find.Execute("Questions: Introduction to complex numbers", $false, $false, $false, $false, $false, `
    $true, 1, $false, "Questions: Introduction to complex numbers", 2)

# Author paragraph: "Tom Coleman"
$authorRange = $d.Paragraphs.Item(2).Range
$authorRange.Find.Execute("Tom Coleman", $false, $false, $false, $false, $false, `
    $true, 1, $false, "Tom Coleman", 2)

# Abstract paragraph: "A selection of questions for the study guide on introduction to complex numbers."
$abstractRange = $d.Paragraphs.Item(4).Range
$abstractRange.Find.Execute("A selection of questions for the study guide on introduction to complex numbers.", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "A selection of questions for the study guide on introduction to complex numbers.", 2)
